$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert the new "ng2-file-upload" row under the Third party libraries section (new row 5),
#    shifting everything below down by one.
$ws.Rows(5).Insert()

# 2. Insert a blank spacer row after the "Live Share" extension row (new row 22 once the
#    previous insert has settled), shifting the DotNet Commands section (and below) down by one more.
$ws.Rows(22).Insert()

# 3. Give the "Third party libraries" header the same bold style used by the other section headers.
$ws.Range("A7").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# 4. Fix the swapped name/description columns for ngx-bootstrap.
$ws.Range("B4").Value = "ngx-bootstrap"
$ws.Range("C4").Value = "The best way to quickly integrate Bootstrap 3 or Bootstrap 4 Components with Angular"

# 5. Populate the new ng2-file-upload row.
$ws.Range("B5").Value = "ng2-file-upload"
$ws.Range("C5").Value = "Easy to use Angular2 directives for files upload (demo)"

# 6. Widen column B to fit the new content and drop the old bestFit auto-sizing.
$ws.Columns(2).ColumnWidth = 59.67

# 7. Move the active selection as it ended up after the edits.
$ws.Range("B4").Select()
